$wb = $excel.ActiveWorkbook

$ovw  = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared by Overview!E2/F2/E3/F3 and the Status column on both language
#    sheets)
# ---------------------------------------------------------------------------
$statusText = "Handed back: in sync with en-US"

$ovw.Range("E2").Value = $statusText
$ovw.Range("F2").Value = $statusText
$ovw.Range("E3").Value = $statusText
$ovw.Range("F3").Value = $statusText

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime:
#    zh-cn keeps using the same (now updated) shared string -> set it once,
#    which will retroactively update every cell that still points at it.
#    de-de gets its own, later timestamp.
# ---------------------------------------------------------------------------
$zhcn.Range("L2").Value = "2016-12-08 07:28:19"
$zhcn.Range("L3").Value = "2016-12-08 07:28:19"

$dede.Range("L2").Value = "2016-12-08 07:28:38"
$dede.Range("L3").Value = "2016-12-08 07:28:38"

# ---------------------------------------------------------------------------
# 3. Populate "Latest Target File" (J) and "Latest Handback File" (K) columns
#    with the generated handback report for both rows on both language
#    sheets.  J gets a hyperlink identical to the one already used in column
#    A (same target file, same display text).
# ---------------------------------------------------------------------------

# zh-cn, row 2 (813e9e03 file)
$zhcn.Hyperlinks.Add($zhcn.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34449e3a5f705fa34b3efbe3271624d9ac1327fa/e2e/813e9e03-4fc6-4ec5-a03f-ecb00fad1373.md", "", "", "813e9e03-4fc6-4ec5-a03f-ecb00fad1373.md") | Out-Null
$zhcn.Range("K2").Value = "813e9e03-4fc6-4ec5-a03f-ecb00fad1373.7e2ab1bc1ad91d30c3f652697fa4868dc732b96b.zh-cn.xlf"

# zh-cn, row 3 (cc64cc7e file)
$zhcn.Hyperlinks.Add($zhcn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34449e3a5f705fa34b3efbe3271624d9ac1327fa/e2e/cc64cc7e-d456-48bd-ad0b-cd1b294a0d16.md", "", "", "cc64cc7e-d456-48bd-ad0b-cd1b294a0d16.md") | Out-Null
$zhcn.Range("K3").Value = "cc64cc7e-d456-48bd-ad0b-cd1b294a0d16.b4b7be6a9a7c30c9feea9ff2086dc878938f0f80.zh-cn.xlf"

# de-de, row 2 (813e9e03 file)
$dede.Hyperlinks.Add($dede.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34449e3a5f705fa34b3efbe3271624d9ac1327fa/e2e/813e9e03-4fc6-4ec5-a03f-ecb00fad1373.md", "", "", "813e9e03-4fc6-4ec5-a03f-ecb00fad1373.md") | Out-Null
$dede.Range("K2").Value = "813e9e03-4fc6-4ec5-a03f-ecb00fad1373.7e2ab1bc1ad91d30c3f652697fa4868dc732b96b.de-de.xlf"

# de-de, row 3 (cc64cc7e file)
$dede.Hyperlinks.Add($dede.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34449e3a5f705fa34b3efbe3271624d9ac1327fa/e2e/cc64cc7e-d456-48bd-ad0b-cd1b294a0d16.md", "", "", "cc64cc7e-d456-48bd-ad0b-cd1b294a0d16.md") | Out-Null
$dede.Range("K3").Value = "cc64cc7e-d456-48bd-ad0b-cd1b294a0d16.b4b7be6a9a7c30c9feea9ff2086dc878938f0f80.de-de.xlf"

# ---------------------------------------------------------------------------
# 4. Column widths widened to fit the longer status text / new file-name
#    columns.
# ---------------------------------------------------------------------------
$ovw.Columns.Item(5).ColumnWidth = 29.2
$ovw.Columns.Item(6).ColumnWidth = 29.2

$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(10).ColumnWidth = 39.15
$zhcn.Columns.Item(11).ColumnWidth = 39.15

$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(10).ColumnWidth = 39.15
$dede.Columns.Item(11).ColumnWidth = 39.15
